# Add 2022-Q1 data:
#  - the current "总计" (Total) sheet is renamed to "2022-Q1" and filled
#    with the fund-level holdings for that quarter (it keeps sheetId/rId
#    of the old "总计" sheet)
#  - a brand-new "总计" sheet is appended right after it (a copy of the
#    original "总计" sheet, so it keeps the same sheet-level formatting),
#    and a new summary row for 2022-Q1 is written into it.

$wb = $excel.ActiveWorkbook

$oldTotal = $wb.Worksheets.Item("总计")

# Make a copy of the existing "总计" sheet; it is placed immediately
# after the original and will become the new, updated "总计" sheet.
$oldTotal.Copy($null, $oldTotal)
$newTotal = $wb.Worksheets.Item("总计 (2)")
$newTotal.Name = "总计temp"

# Repurpose the original "总计" sheet as the "2022-Q1" fund detail sheet.
$oldTotal.Name = "2022-Q1"
$q1 = $oldTotal

# --- Build "2022-Q1" fund detail sheet -------------------------------
$template = $wb.Worksheets.Item("2021-Q4")
$template.Range("A1:H13").Copy()
$q1.Range("A1").PasteSpecial(-4122)   # xlPasteFormats

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Columns B, D, E, F, G hold numeric-looking text (fund codes with
# leading zeros, percentages, etc.) that must stay text.
$q1.Range("B2:B13").NumberFormat = "@"
$q1.Range("D2:G13").NumberFormat = "@"

$funds = @(
  @("516150", "嘉实中证稀土产业ETF", "25.17", "99.75", "5.42", "1.3642", 4),
  @("516780", "华泰柏瑞中证稀土产业ETF", "11.06", "98.70", "5.32", "0.5884", 4),
  @("159715", "易方达中证稀土产业ETF", "3.42", "99.06", "5.38", "0.1840", 4),
  @("159713", "富国中证稀土产业交易型开放式指数证券投资基金", "3.26", "99.26", "5.41", "0.1764", 4),
  @("159608", "广发中证稀有金属ETF", "2.39", "98.72", "2.58", "0.0617", 10),
  @("159962", "华夏中证四川国企改革ETF", "0.49", "95.82", "4.83", "0.0237", 5),
  @("014331", "华泰柏瑞中证稀土产业ETF联接A", "0.86", "24.22", "1.32", "0.0114", 6),
  @("014332", "华泰柏瑞中证稀土产业ETF联接C", "0.70", "24.22", "1.32", "0.0092", 6),
  @("003854", "汇安丰华灵活配置混合A", "0.19", "29.16", "1.29", "0.0025", 10),
  @("003855", "汇安丰华灵活配置混合C", "0.19", "29.16", "1.29", "0.0025", 10),
  @("006601", "国融融泰灵活配置混合A", "0.04", "47.44", "4.12", "0.0016", 2),
  @("006602", "国融融泰灵活配置混合C", "0.01", "47.44", "4.12", "0.0004", 2)
)

for ($i = 0; $i -lt $funds.Count; $i++) {
  $r = $i + 2
  $row = $funds[$i]
  $q1.Range("A$r").Value = $i
  $q1.Range("B$r").Value = $row[0]
  $q1.Range("C$r").Value = $row[1]
  $q1.Range("D$r").Value = $row[2]
  $q1.Range("E$r").Value = $row[3]
  $q1.Range("F$r").Value = $row[4]
  $q1.Range("G$r").Value = $row[5]
  $q1.Range("H$r").Value = $row[6]
}

# --- Build the updated "总计" summary sheet ---------------------------
# Add one more row (copying formatting from the existing last data row)
# to make room for the new 2022-Q1 summary entry.
$newTotal.Range("A6:D6").Copy()
$newTotal.Range("A7").PasteSpecial(-4122)   # xlPasteFormats

$summary = @(
  @(0, "2022-Q1", 12, 2.43),
  @(1, "2021-Q4", 13, 3.87),
  @(2, "2021-Q3", 9, 3.08),
  @(3, "2021-Q2", 20, 2.99),
  @(4, "2021-Q1", 17, 3.28),
  @(5, "2020-Q4", 4, 0.17)
)

for ($i = 0; $i -lt $summary.Count; $i++) {
  $r = $i + 2
  $row = $summary[$i]
  $newTotal.Range("A$r").Value = $row[0]
  $newTotal.Range("B$r").Value = $row[1]
  $newTotal.Range("C$r").Value = $row[2]
  $newTotal.Range("D$r").Value = $row[3]
}

$newTotal.Name = "总计"

# Leave the first sheet active (clears the "tab selected" flag that
# landed on the copied sheet).
$wb.Worksheets.Item(1).Activate()
